$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'42.522.59"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = "'  +0.06%  "
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(3, 4).Value = "'2.520.31"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = "'  -1.09%  "
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(4, 4).Value = "'0.998"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = "'  -0.04%  "
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(5, 4).Value = "'313.68"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = "'  +0.34%  "
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(6, 4).Value = "'99.00"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = "'  -1.65%  "
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = "'  -1.34%  "
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = "'  -0.02%  "
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = "'  -2.56%  "
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(10, 4).Value = "'35.16"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = "'  -2.50%  "
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(11, 4).Value = "'0.0799"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = "'  -0.61%  "
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = "'  +0.92%  "
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = "'  -2.89%  "
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(14, 4).Value = "'2.901.21"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = "'  -1.07%  "
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(15, 4).Value = "'15.27"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = "'  -4.72%  "
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(16, 4).Value = "'2.518.71"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = "'  -0.89%  "
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(17, 4).Value = "'0.805"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = "'  -4.19%  "
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(18, 4).Value = "'42.484.06"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = "'  -0.12%  "
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(19, 4).Value = "'6.60"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = "'  -2.91%  "
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(20, 4).Value = "'0.0₃0938"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = "'  -1.54%  "
$ws.Cells.Item(20, 5).Style = 'Normal'
$ws.Cells.Item(21, 4).Value = "'12.08"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = "'  -1.96%  "
$ws.Cells.Item(21, 5).Style = 'Normal'
$ws.Cells.Item(22, 4).Value = "'68.97"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = "'  -0.25%  "
$ws.Cells.Item(22, 5).Style = 'Normal'
$ws.Cells.Item(23, 4).Value = "'240.91"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = "'  -1.42%  "
$ws.Cells.Item(23, 5).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = "'  -2.21%  "
$ws.Cells.Item(24, 5).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = "'  -3.62%  "
$ws.Cells.Item(25, 5).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = "'  +0.03%  "
$ws.Cells.Item(26, 5).Style = 'Normal'
$ws.Cells.Item(27, 4).Value = "'25.36"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = "'  -3.85%  "
$ws.Cells.Item(27, 5).Style = 'Normal'
$ws.Cells.Item(28, 4).Value = "'2.25"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = "'  -4.66%  "
$ws.Cells.Item(28, 5).Style = 'Normal'
$ws.Cells.Item(29, 4).Value = "'9.98"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = "'  -1.34%  "
$ws.Cells.Item(29, 5).Style = 'Normal'
$ws.Cells.Item(30, 4).Value = "'37.66"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = "'  -7.01%  "
$ws.Cells.Item(30, 5).Style = 'Normal'
$ws.Cells.Item(31, 4).Value = "'5.88"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = "'  +3.05%  "
$ws.Cells.Item(31, 5).Style = 'Normal'
$ws.Cells.Item(32, 4).Value = "'156.80"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = "'  -1.04%  "
$ws.Cells.Item(32, 5).Style = 'Normal'
$ws.Cells.Item(33, 4).Value = "'2.68"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = "'  -4.23%  "
$ws.Cells.Item(33, 5).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = "'  +0.69%  "
$ws.Cells.Item(34, 5).Style = 'Normal'
$ws.Cells.Item(35, 4).Value = "'0.0781"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = "'  -2.75%  "
$ws.Cells.Item(35, 5).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = "'  -2.00%  "
$ws.Cells.Item(36, 5).Style = 'Normal'
$ws.Cells.Item(37, 2).Value = "'Celestia"
$ws.Cells.Item(37, 2).Style = 'Normal'
$ws.Cells.Item(37, 3).Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(37, 3).Style = 'Normal'
$ws.Cells.Item(37, 4).Value = "'17.73"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = "'  -2.32%  "
$ws.Cells.Item(37, 5).Style = 'Normal'
$ws.Cells.Item(38, 2).Value = "'ARBITRUM"
$ws.Cells.Item(38, 2).Style = 'Normal'
$ws.Cells.Item(38, 3).Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(38, 3).Style = 'Normal'
$ws.Cells.Item(38, 4).Value = "'1.96"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = "'  -4.97%  "
$ws.Cells.Item(38, 5).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = "'  -3.34%  "
$ws.Cells.Item(39, 5).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = "'  -0.83%  "
$ws.Cells.Item(40, 5).Style = 'Normal'
$ws.Cells.Item(41, 4).Value = "'4.15"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = "'  -1.12%  "
$ws.Cells.Item(41, 5).Style = 'Normal'
$ws.Cells.Item(42, 4).Value = "'22.02"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = "'  -0.50%  "
$ws.Cells.Item(42, 5).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = "'  -0.12%  "
$ws.Cells.Item(43, 5).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = "'  -1.17%  "
$ws.Cells.Item(44, 5).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = "'  -3.74%  "
$ws.Cells.Item(45, 5).Style = 'Normal'
$ws.Cells.Item(46, 4).Value = "'1.994.15"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = "'  +1.30%  "
$ws.Cells.Item(46, 5).Style = 'Normal'
$ws.Cells.Item(47, 4).Value = "'8.95"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = "'  +0.81%  "
$ws.Cells.Item(47, 5).Style = 'Normal'
$ws.Cells.Item(48, 4).Value = "'2.756.03"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = "'  -1.15%  "
$ws.Cells.Item(48, 5).Style = 'Normal'
$ws.Cells.Item(49, 4).Value = "'78.76"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = "'  -2.88%  "
$ws.Cells.Item(49, 5).Style = 'Normal'
$ws.Cells.Item(50, 4).Value = "'0.187"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = "'  -3.02%  "
$ws.Cells.Item(50, 5).Style = 'Normal'
$ws.Cells.Item(51, 4).Value = "'71.38"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = "'  -2.44%  "
$ws.Cells.Item(51, 5).Style = 'Normal'
